$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 2, 4, 5, 6, 8 per repulled data
$ws.Range("F2").Value = -4
$ws.Range("F4").Value = 3
$ws.Range("F5").Value = -2
$ws.Range("F6").Value = 4
$ws.Range("F8").Value = -2
